# Update cryptocurrency price and volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay text so values are not
# reinterpreted as numbers/dates by Excel (e.g. "1.00", "51.891.80").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "51.891.80"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "2.922.64"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "361.02"
$ws.Range("E5").Value = "  +2.30%  "
$ws.Range("D6").Value = "110.41"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("D9").Value = "0.633"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "39.38"
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("D11").Value = "0.0880"
$ws.Range("E11").Value = "  +1.75%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").Value = "7.88"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "3.388.32"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "2.911.01"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "0.991"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "51.890.08"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "3.36"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "7.62"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").Value = "14.11"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "71.14"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "270.31"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "2.86"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  +13.28%  "
$ws.Range("D27").Value = "27.05"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").Value = "7.67"
$ws.Range("E28").Value = "  +17.45%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +13.28%  "
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "38.27"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "2.28"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").Value = "52.43"
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("D39").Value = "18.46"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("E40").Value = "  -3.76%  "
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("D43").Value = "23.04"
$ws.Range("E43").Value = "  -4.99%  "
$ws.Range("D44").Value = "119.46"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").Value = "3.51"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").Value = "2.129.98"
$ws.Range("E47").Value = "  -4.17%  "
$ws.Range("E48").Value = "  -5.53%  "
$ws.Range("D49").Value = "0.0336"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").Value = "9.18"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  -5.38%  "
